# Update crypto price/volume data in the worksheet to reflect the
# refreshed values from the "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Setting NumberFormat to "@" (Text) before assigning the value ensures
# Excel keeps the exact textual representation (e.g. "12.00", "0.0000231",
# "61.970.87") instead of silently re-interpreting it as a number and
# dropping trailing zeros / using scientific notation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.970.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.904.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.70"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.901.58"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.04"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.387.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.933.79"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.899.81"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.85"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.16"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -9.66%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000112"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.52"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.61"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.38"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.76"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.32%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.14%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.11"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.14"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.703.45"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "346.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.20%  "
